$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Programs/Studies/Cases/Samples/Files" query - used in column C for rows 2,3,4
$programsQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed  IN ['Saint Bernard']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Cases query (row 2, column B) - appended Cohort column to the RETURN clause
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Saint Bernard']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# Files query (row 4, column B) - removed trailing Study Code line
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Saint Bernard']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
          coalesce(diag.disease_term,'') AS Diagnosis 
'@

function Trim-TrailingNewline([string]$s) {
    if ($s.EndsWith("`r`n")) { return $s.Substring(0, $s.Length - 2) }
    if ($s.EndsWith("`n")) { return $s.Substring(0, $s.Length - 1) }
    return $s
}

$programsQuery = Trim-TrailingNewline($programsQuery)
$casesQuery = Trim-TrailingNewline($casesQuery)
$filesQuery = Trim-TrailingNewline($filesQuery)

# Write order controls the resulting shared-string table order (matching the
# target file's layout): new Programs query first, then the modified Files
# query, then the modified Cases query. B3 (sample query) is left untouched.
$ws.Cells.Item(2, 3).Value = $programsQuery
$ws.Cells.Item(3, 3).Value = $programsQuery
$ws.Cells.Item(4, 3).Value = $programsQuery

$ws.Cells.Item(4, 2).Value = $filesQuery

$ws.Cells.Item(2, 2).Value = $casesQuery

# Row heights shrink to fit the shorter query text (values from the target file)
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# Selection moves to B2 per the edited file
$ws.Range("B2").Select()
